$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.353.77'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.20%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.847.36'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.14%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9988'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6271'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.38%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9994'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07602'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.18%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2901'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.34%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.67'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.51%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07740'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.11%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.026'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.04%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6783'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.45%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.00001066'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.39%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.94'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.00%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.127'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.39%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.395.88'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.13%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '227.73'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.65%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.35'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.99%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9995'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.08%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.505'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.72%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9994'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.09%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '158.69'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.99%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1383'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.18%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.433'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.55%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '17.64'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.18%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.431'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.73%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.469'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.07%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05619'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.33%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.099'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.39%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.064'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.36%  '

$ws.Range("B32").Value = 'ARBITRUM'
$ws.Range("C32").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.160'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.03%  '

$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.829'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.6973'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.38%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.583'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.19%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.231.36'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.52%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01798'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.10%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.723'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.364'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9000'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.95%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9990'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.13%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '101.49'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.08%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.40'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.201'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.88%  '

$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3989'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.89%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.000'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.30%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.683'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.61%  '

$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1140'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.26%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05697'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.24%  '

$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4625'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.04%  '

$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.342'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.56%  '
